# Zeitmanagement.xlsx - add "Total" sheet summarizing planned/actual hours
# across the four per-person sheets (Agdas, Hussein, Sentler, Brak, Kessener).

$wb = $excel.ActiveWorkbook

# --- minor view-state change observed on the "Sentler" sheet (selection moved) ---
$sentler = $wb.Worksheets.Item("Sentler")
$sentler.Range("C2").Select()

# --- add the new "Total" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Total"

# Row 1: planned time total
$newSheet.Range("A1").Value = "Geplante Zeit:"
$newSheet.Range("C1").Formula = "=SUM(Agdas!C2:C100,Hussein!C2:C100,Sentler!C2:C100,Brak!C2:C100,Kessener!C2:C100)"

# Row 2: actual ("tatsaechliche") time total
$newSheet.Range("A2").Value = "Tatsächsliche Zeit:"
$newSheet.Range("C2").Formula = "=SUM(Agdas!D2:D100,Hussein!D2:D100,Sentler!D2:D100,Brak!D2:D100,Kessener!D2:D100)"
$newSheet.Range("C2").NumberFormat = "# ?/?"

# Right-align + merge the label cells across A:B for both rows
$newSheet.Range("A1:B2").HorizontalAlignment = -4152
$newSheet.Range("A1:B1").Merge()
$newSheet.Range("A2:B2").Merge()

# Page setup matching the rest of the workbook
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Leave selection on E5 / Total tab active, matching final saved view state
$newSheet.Range("E5").Select()
